$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.625.35'
$ws.Range("E2").Value = '  -1.34%  '

$ws.Range("D3").Value = '1.631.90'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'211.25"
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("E6").Value = '  -1.01%  '

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").Value = "'23.03"
$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("E9").Value = '  -0.23%  '

$ws.Range("D10").Value = "'0.0611"
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D11").Value = "'0.0862"
$ws.Range("E11").Value = '  -3.34%  '

$ws.Range("D12").Value = '1.860.18'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").Value = '1.628.73'
$ws.Range("E13").Value = '  -0.92%  '

$ws.Range("E14").Value = '  -0.52%  '

$ws.Range("D15").Value = "'0.558"
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = "'65.10"
$ws.Range("E16").Value = '  +0.54%  '

$ws.Range("D17").Value = '27.588.56'
$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").Value = "'229.68"
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("D19").Value = '0.0₃0719'
$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = '  -1.35%  '

$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").Value = "'10.64"
$ws.Range("E22").Value = '  +6.14%  '

$ws.Range("D23").Value = "'4.37"
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("E24").Value = '  +2.68%  '

$ws.Range("D25").Value = "'149.06"
$ws.Range("E25").Value = '  -1.06%  '

$ws.Range("D26").Value = "'6.87"
$ws.Range("E26").Value = '  -1.14%  '

$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").Value = "'15.61"
$ws.Range("E28").Value = '  -0.64%  '

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  -0.60%  '

$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").Value = '1.463.61'
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = '  -0.79%  '

$ws.Range("E36").Value = '  -1.72%  '

$ws.Range("D37").Value = "'0.880"
$ws.Range("E37").Value = '  -0.42%  '

$ws.Range("D38").Value = "'0.558"
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = "'68.77"
$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").Value = "'2.45"
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = '  -0.94%  '

$ws.Range("D46").Value = "'5.37"
$ws.Range("E46").Value = '  -0.93%  '

$ws.Range("D47").Value = '1.771.73'
$ws.Range("E47").Value = '  -0.77%  '

$ws.Range("E48").Value = '  +1.60%  '

$ws.Range("D49").Value = "'87.42"
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("E50").Value = '  -1.05%  '

$ws.Range("D51").Value = "'0.0996"
$ws.Range("E51").Value = '  +0.07%  '

# Reset quote-prefix styling introduced by forcing text on numeric-looking values
$textCells = @("D5","D7","D8","D10","D11","D15","D16","D18","D20","D21","D22","D23","D25","D26","D28","D29","D35","D37","D38","D41","D42","D43","D44","D45","D46","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
